# Weekly fruit/vegetable update: a new daily price record was inserted
# above the existing row 199 (Feria Lagunitas de Puerto Montt - Perejil),
# pushing the previous rows 199-212 down to 200-213.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 199 (shifts 199:212 down to 200:213)
$ws.Rows.Item(199).Insert()

# Populate the newly inserted row with the new record's data
$ws.Range("A199").Value = 4
$ws.Range("B199").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C199").Value = "Los Lagos"
$ws.Range("D199").Value = 44610
$ws.Range("E199").Value = 10
$ws.Range("F199").Value = 100112044
$ws.Range("G199").Value = "Perejil"
$ws.Range("H199").Value = "Sin especificar"
$ws.Range("I199").Value = "Primera"
$ws.Range("J199").Value = 150
$ws.Range("K199").Value = 5000
$ws.Range("L199").Value = 5000
$ws.Range("M199").Value = 5000
$ws.Range("N199").Value = "$/docena de atados (3 kilos)"
$ws.Range("O199").Value = "Región Metropolitana"
$ws.Range("P199").Value = 1667
$ws.Range("Q199").Value = 3
$ws.Range("R199").Value = "Hortaliza"
